$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) / volume-change (E) figures per row, plus the
# Monero / EthereumClassic row swap (rows 38 and 39).
$data = @{
  2  = @{ D = "68.742.56";  E = "  +0.73%  " }
  3  = @{ D = "2.710.48";   E = "  +2.39%  " }
  4  = @{ D = "0.999";      E = "  -0.09%  " }
  5  = @{ D = "600.00";     E = "  +0.42%  " }
  6  = @{ D = "162.65";     E = "  +3.82%  " }
  7  = @{               E = "  -0.02%  " }
  8  = @{               E = "  +0.21%  " }
  9  = @{ D = "2.708.69";   E = "  +2.35%  " }
  10 = @{               E = "  +0.36%  " }
  11 = @{               E = "  -0.31%  " }
  12 = @{ D = "5.32";       E = "  +1.13%  " }
  13 = @{ D = "0.362";      E = "  +3.15%  " }
  14 = @{ D = "28.47";      E = "  +1.65%  " }
  15 = @{ D = "3.194.05";   E = "  +2.08%  " }
  16 = @{               E = "  -0.54%  " }
  17 = @{ D = "68.651.96";  E = "  +0.69%  " }
  18 = @{ D = "2.731.35";   E = "  +3.23%  " }
  19 = @{               E = "  +4.46%  " }
  20 = @{ D = "7.68";       E = "  +4.69%  " }
  21 = @{ D = "365.55";     E = "  +0.58%  " }
  22 = @{               E = "  +3.00%  " }
  23 = @{ D = "4.94";       E = "  +2.87%  " }
  24 = @{               E = "  +2.53%  " }
  25 = @{ D = "74.19";      E = "  -1.30%  " }
  26 = @{               E = "  -0.02%  " }
  27 = @{ D = "9.92";       E = "  +1.64%  " }
  28 = @{               E = "  +2.13%  " }
  29 = @{               E = "  +1.37%  " }
  30 = @{ D = "598.44";     E = "  +6.85%  " }
  31 = @{ D = "1.00";       E = "  +0.07%  " }
  32 = @{ D = "8.30";       E = "  +3.07%  " }
  33 = @{               E = "  +3.05%  " }
  34 = @{               E = "  +4.84%  " }
  35 = @{ D = "0.133";      E = "  +3.50%  " }
  36 = @{               E = "  +5.52%  " }
  40 = @{ D = "0.380";      E = "  +2.44%  " }
  41 = @{               E = "  +2.55%  " }
  42 = @{ D = "5.44";       E = "  +2.18%  " }
  43 = @{ D = "2.72";       E = "  +4.27%  " }
  44 = @{               E = "  +1.22%  " }
  45 = @{               E = "  -5.25%  " }
  46 = @{               E = "  +0.04%  " }
  47 = @{ D = "158.43";     E = "  -0.29%  " }
  48 = @{ D = "3.95";       E = "  +5.97%  " }
  49 = @{               E = "  +5.60%  " }
  50 = @{               E = "  +7.25%  " }
  51 = @{ D = "22.16";      E = "  +0.42%  " }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
    }
}

# Rows 38 and 39 swap places: Monero <-> EthereumClassic (name, link,
# price and volume change all move together).
$ws.Cells.Item(38, 2).Value = "EthereumClassic"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "19.88"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  +1.01%  "

$ws.Cells.Item(39, 2).Value = "Monero"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "160.34"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  -0.94%  "
